{"js": "// The June sitrep R output dump (two \"SourceCode\" paragraphs showing the\n// `## # A tibble: ...` console output) is being removed, and the paragraph\n// that follows them (\"UNICEF also supported the training of 148 ...\")\n// reverts from the \"First Paragraph\" style to plain \"Body Text\" (since it's\n// no longer the first paragraph following a heading).\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text,items/style\");\nawait context.sync();\n\n// Locate the two verbatim \"tibble\" dump paragraphs by their distinctive text.\nconst toDelete = [];\nlet lastDeletedIndex = -1;\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  const p = paragraphs.items[i];\n  const text = p.text || \"\";\n  if (text.indexOf(\"## # A tibble: 0 x 10\") !== -1 || text.indexOf(\"## # A tibble: 17,052 x 10\") !== -1) {\n    toDelete.push(p);\n    lastDeletedIndex = i;\n  }\n}\n\n// The paragraph immediately after the last deleted \"tibble\" dump is the one\n// whose style reverts from \"First Paragraph\" to \"Body Text\".\nconst nextParagraph = lastDeletedIndex >= 0 ? paragraphs.items[lastDeletedIndex + 1] : null;\nif (nextParagraph) {\n  nextParagraph.style = \"Body Text\";\n}\n\nfor (const p of toDelete) {\n  p.delete();\n}\n\nawait context.sync();\n", "ps1": "# The June sitrep R output dump (two \"Source Code\" styled paragraphs showing\n# the `## # A tibble: ...` console output) is being removed, and the\n# paragraph that follows them (\"UNICEF also supported the training of 148\n# ...\") reverts from the \"First Paragraph\" style to plain \"Body Text\" (since\n# it's no longer the first paragraph following a heading).\n\n$d = $word.ActiveDocument\n\n# Locate the two verbatim \"tibble\" dump paragraphs by their distinctive text.\n$indices = @()\nfor ($i = 1; $i -le $d.Paragraphs.Count; $i++) {\n    $t = $d.Paragraphs.Item($i).Range.Text\n    if ($t -match [regex]::Escape(\"## # A tibble: 0 x 10\") -or $t -match [regex]::Escape(\"## # A tibble: 17,052 x 10\")) {\n        $indices += $i\n    }\n}\n\n# The paragraph immediately after the last deleted \"tibble\" dump is the one\n# whose style reverts from \"First Paragraph\" to \"Body Text\".\n$lastIndex = ($indices | Measure-Object -Maximum).Maximum\nif ($lastIndex) {\n    $nextParagraph = $d.Paragraphs.Item($lastIndex + 1)\n    $nextParagraph.Range.Style = \"Body Text\"\n}\n\n# Delete from the bottom up so earlier paragraph indices stay valid.\n$sortedDescending = $indices | Sort-Object -Descending\nforeach ($index in $sortedDescending) {\n    $d.Paragraphs.Item($index).Range.Delete()\n}\n"}
